# Insert one new data row at row 551 (pushes existing rows 551-616 down to
# 552-617) and populate it with the new "Femacal de La Calera - Ajo" record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(551).Insert()

$ws.Range("A551").Value = 3
$ws.Range("B551").Value = "Femacal de La Calera"
$ws.Range("C551").Value = "Coquimbo"
$ws.Range("D551").Value = 44918
$ws.Range("E551").Value = 5
$ws.Range("F551").Value = 100112003
$ws.Range("G551").Value = "Ajo"
$ws.Range("H551").Value = "Chino"
$ws.Range("I551").Value = "Primera"
$ws.Range("J551").Value = 75
$ws.Range("K551").Value = 13000
$ws.Range("L551").Value = 13500
$ws.Range("M551").Value = 13233
$ws.Range("N551").Value = '$/caja 10 kilos'
$ws.Range("O551").Value = "China"
$ws.Range("P551").Value = 1323
$ws.Range("Q551").Value = 10
$ws.Range("R551").Value = "Hortaliza"

# Keep the date column's existing date-style formatting consistent with the
# rest of column D.
$ws.Range("D551").NumberFormat = $ws.Range("D552").NumberFormat
